# test code commit by Kevin
#
# Row 2 on both template sheets holds one sample/demo record. This refreshes
# that demo record with a new invoice number, buyer/supplier reference and
# contact names (the PO sheet mirrors invoice/PO numbers in swapped
# columns relative to the Invoice sheet).

$wb = $excel.ActiveWorkbook

$wsInvoice = $wb.Worksheets.Item("Historical Invoice Template")
$wsPO      = $wb.Worksheets.Item("Historical PO Template")

# Historical Invoice Template, row 2: Invoice No. / PO No. / Supplier / Buyer
$wsInvoice.Range("B2").Value = "Invoice1393659"
$wsInvoice.Range("C2").Value = "Invoice1881222"
$wsInvoice.Range("E2").Value = "Laura036q"
$wsInvoice.Range("F2").Value = "Lucye7o4"

# Historical PO Template, row 2: PO No. / Invoice No. / Supplier / Buyer
$wsPO.Range("B2").Value = "Invoice1881222"
$wsPO.Range("C2").Value = "Invoice1393659"
$wsPO.Range("E2").Value = "Laura036q"
$wsPO.Range("F2").Value = "Lucye7o4"
